# Logged Week 15 and simulated Week 16
# Update the "Road" (R) row target-depth splits on both the OFF and DEF
# sheets to reflect the newly logged/simulated week totals.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: Short Att / Short Comp / Deep Att / Deep Comp for "R" row ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 357
$wsOff.Range("C3").Value = 250
$wsOff.Range("D3").Value = 82
$wsOff.Range("E3").Value = 33

# --- DEF sheet: Short Att / Short Comp / Deep Att / Deep Comp for "R" row ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 370
$wsDef.Range("C3").Value = 255
$wsDef.Range("D3").Value = 116
$wsDef.Range("E3").Value = 60
